$d = $word.ActiveDocument

function Add-RunAfterParagraph($paragraph, $text) {
    # Appends $text to the end of $paragraph's text as a brand-new, separate
    # <w:r> run (rather than being merged into the paragraph's final existing
    # run). This is done by temporarily splitting the paragraph in two,
    # typing the new text into the freshly-created (and therefore
    # independent) paragraph, and then deleting the paragraph mark that
    # separates them again -- which re-joins the two paragraphs without
    # forcing a reformat/merge of their runs.
    $r = $paragraph.Range
    $endPos = $r.End - 1
    $splitPoint = $d.Range($endPos, $endPos)
    $splitPoint.InsertParagraphAfter()
    $newPara = $paragraph.Next()
    $r2 = $newPara.Range
    $insertPoint = $d.Range($r2.Start, $r2.Start)
    $insertPoint.InsertAfter($text)
    $markRange = $d.Range($endPos, $r2.Start)
    $markRange.Delete()
}

function Merge-ParagraphRuns($paragraph, $text) {
    # Collapses every run in $paragraph's text into a single run containing
    # $text, using Find/Replace over the paragraph range (excluding the
    # trailing paragraph mark).
    $r = $paragraph.Range
    $scoped = $d.Range($r.Start, $r.End - 1)
    $scoped.Find.Execute($scoped.Text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# 1. "Estimated range to target" -> append ", in meters" as a new run
Add-RunAfterParagraph $d.Paragraphs(17) ", in meters"

# 2. Merge "Estimated " / "Doppler range rate of" / " target" into a single
#    run, then append ", in meters per second" as a new run
Merge-ParagraphRuns $d.Paragraphs(20) "Estimated Doppler range rate of target"
Add-RunAfterParagraph $d.Paragraphs(20) ", in meters per second"

# 3. Merge "Estimated " / "azimuth bearing of target" into a single run,
#    then append ", in degrees" as a new run
Merge-ParagraphRuns $d.Paragraphs(23) "Estimated azimuth bearing of target"
Add-RunAfterParagraph $d.Paragraphs(23) ", in degrees"

# 4. Merge "'" / "steer" / "'" into a single run "'steer'"
Merge-ParagraphRuns $d.Paragraphs(25) "‘steer’"

# 5. "Antenna array steering direction" -> append ", in degrees" as a new run
Add-RunAfterParagraph $d.Paragraphs(26) ", in degrees"

# 6. Merge "'" / "SNR" / "'" into a single run "'SNR'"
Merge-ParagraphRuns $d.Paragraphs(28) "‘SNR’"

# 7. Merge "Estimated " / "signal-to-noise ratio of target reflection" into a
#    single run, then append ", in decibels" as a new run
Merge-ParagraphRuns $d.Paragraphs(29) "Estimated signal-to-noise ratio of target reflection"
Add-RunAfterParagraph $d.Paragraphs(29) ", in decibels"

# 8. Merge "Boolean " / "(# of radar units) x (# of frames) array" into a
#    single run
Merge-ParagraphRuns $d.Paragraphs(33) "Boolean (# of radar units) x (# of frames) array"
